$d = $word.ActiveDocument

# Locate the anchor paragraph (the last populated "Paragrafo elenco" item before the
# block of trailing empty paragraphs at the end of the document), then target the
# empty paragraph that immediately follows it.
$r = $d.Content
$r.Find.ClearFormatting()
$ok = $r.Find.Execute("l’applicazione.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Anchor text not found"
}
$pos = $r.Start

$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $pos -and $p.Range.End -ge $pos) {
        $targetIndex = $i + 1
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Mettere un controllo in più sull’Eseguito: scriverlo solamente se effettivamente sono stati copiati i dati di quel giorno</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Considerare importi negativi per INCASSI e PROVVIGIONI</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Per le PROVVIGIONI devo considerare i MOBILE POS</w:t></w:r><w:r><w:t xml:space="preserve"> e VIRTUAL POS</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Ricontrollare se considero FINANZIAMENTO AL CONSUMO</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Chiedere a GIGI di riscaricare il file di CATTOLICA del 01/03</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Inserire errore se metodo di pagamento nuovo</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Togliere errore se da un file di GENERALI, CATTOLICA o TUTELA LEGALE sono vuoti quindi non si legge la data: Data mancante</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Chiedere a Gigi un file di CATTOLICA in cui vi è un finanziamento al consumo ‘Bonifico su CC di Direzione’, per vedere se vengono fatte più righe o meno, in modo tale da capire se considerare solo le provvigioni o anche gli incassi per tale metodo di pagamento</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target = $d.Paragraphs.Item($targetIndex).Range
[void]$target.InsertXML($xml)
